$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"  = 5.885
    "A3"  = -21.573
    "D3"  = -7.768000000000001
    "D12" = -7.243
    "A14" = -21.559
    "A21" = -20.469
    "A23" = -20.428
    "D24" = -7.678999999999999
    "A25" = -21.275
    "B25" = 6.575999999999999
    "D25" = -8.146000000000001
    "A26" = -21.351
    "B27" = 6.202
    "A29" = -21.219
    "B31" = 6.657999999999999
    "B39" = 7.517
    "B48" = 5.274
    "D50" = -7.989000000000002
    "B51" = 6.047000000000001
    "B52" = 5.813000000000001
    "A53" = -22.079
    "D53" = -7.704000000000001
    "B55" = 4.671
    "B56" = 5.003
    "A57" = -21.303
    "B57" = 6.431999999999999
    "D57" = -8.231999999999999
    "A59" = -22.397
    "D61" = -7.811
    "D63" = -7.359999999999999
    "A69" = -21.531
    "D70" = -7.163999999999999
    "B73" = 6.751
    "A79" = -21.025
    "A83" = -22.022
    "D86" = -8.07
    "B89" = 5.832000000000001
    "B90" = 5.671
    "A91" = -21.533
    "B92" = 5.899
    "A93" = -21.301
    "D98" = -8.063000000000001
    "D100" = -8.145999999999999
    "D102" = -8.068000000000001
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
